{"js": "// Update three cell descriptions in the \"Atributos de calidad\" table to\n// match the revised wording from the commit.\nconst replacements = [\n  {\n    oldText: \"Posibilidad de a\u00f1adir mayor funcionalidad a\u00f1adiendo nuevos formatos de entrada de datos.\",\n    newText: \"Posibilidad de a\u00f1adir mayor funcionalidad a\u00f1adiendo nuevos formatos de entrada de datos, y de a\u00f1adir nuevas funcionalidad al juego en s\u00ed.\"\n  },\n  {\n    oldText: \"Interfaz sencilla e intuitiva tanto para el operador como para el usuario.\",\n    newText: \"Interfaz sencilla e intuitiva, tanto para el operador de las etapas de extracci\u00f3n, como para el usuario que ejecuta la aplicaci\u00f3n.\"\n  },\n  {\n    oldText: \"Interfaz de juego accesible para todo tipo de personas con diferentes perfiles.\",\n    newText: \"Interfaz de juego accesible para todo tipo de usuarios y dispositivos con diferentes perfiles.\"\n  }\n];\n\nconst body = context.document.body;\n\nfor (const { oldText, newText } of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText);\n  }\n\n  results.items[0].insertText(newText, \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# Update three cell descriptions in the \"Atributos de calidad\" table to\n# match the revised wording from the commit.\n$d = $word.ActiveDocument\n\n$wdFindContinue = 1\n$wdReplaceAll = 2\n\n$replacements = @(\n    @{\n        Old = \"Posibilidad de a\u00f1adir mayor funcionalidad a\u00f1adiendo nuevos formatos de entrada de datos.\"\n        New = \"Posibilidad de a\u00f1adir mayor funcionalidad a\u00f1adiendo nuevos formatos de entrada de datos, y de a\u00f1adir nuevas funcionalidad al juego en s\u00ed.\"\n    },\n    @{\n        Old = \"Interfaz sencilla e intuitiva tanto para el operador como para el usuario.\"\n        New = \"Interfaz sencilla e intuitiva, tanto para el operador de las etapas de extracci\u00f3n, como para el usuario que ejecuta la aplicaci\u00f3n.\"\n    },\n    @{\n        Old = \"Interfaz de juego accesible para todo tipo de personas con diferentes perfiles.\"\n        New = \"Interfaz de juego accesible para todo tipo de usuarios y dispositivos con diferentes perfiles.\"\n    }\n)\n\nforeach ($r in $replacements) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute(\n        $r.Old,\n        $false,\n        $false,\n        $false,\n        $false,\n        $false,\n        $true,\n        $wdFindContinue,\n        $false,\n        $r.New,\n        $wdReplaceAll\n    )\n    if (-not $found) {\n        throw \"Text not found: $($r.Old)\"\n    }\n}\n\n$d.Save()\n"}
